$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refreshed crypto price/volume snapshot (GitHub Actions data pull).
# Most rows only update Price (D) and Volume 1h (E); rows 17-19 also
# reshuffle rank order (Chainlink / Uniswap / WrappedEther).

# Row 2
$ws.Range("D2").Value = "70.840.34"
$ws.Range("E2").Value = "  +2.96%  "

# Row 3
$ws.Range("D3").Value = "3.570.12"
$ws.Range("E3").Value = "  +2.06%  "

# Row 4
$ws.Range("E4").Value = "  -0.01%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "582.06"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.08%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "186.41"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +2.04%  "

# Row 7
$ws.Range("E7").Value = "  +2.16%  "

# Row 8
$ws.Range("D8").Value = "3.557.78"
$ws.Range("E8").Value = "  +1.87%  "

# Row 9
$ws.Range("E9").Value = "  -0.11%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.225"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +23.11%  "

# Row 11
$ws.Range("E11").Value = "  +0.87%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.67"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +1.37%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000321"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +7.00%  "

# Row 14
$ws.Range("E14").Value = "  +0.85%  "

# Row 15
$ws.Range("D15").Value = "4.137.07"
$ws.Range("E15").Value = "  +1.94%  "

# Row 16
$ws.Range("D16").Value = "70.873.83"
$ws.Range("E16").Value = "  +3.19%  "

# Row 17
$ws.Range("B17").Value = "Chainlink"
$ws.Range("C17").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "19.22"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.05%  "

# Row 18
$ws.Range("B18").Value = "Uniswap"
$ws.Range("C18").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.84"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +5.05%  "

# Row 19
$ws.Range("B19").Value = "WrappedEther"
$ws.Range("C19").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D19").Value = "3.568.43"
$ws.Range("E19").Value = "  +2.61%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "573.41"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +5.53%  "

# Row 21
$ws.Range("E21").Value = "  +0.65%  "

# Row 22
$ws.Range("E22").Value = "  -0.67%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.67"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -6.67%  "

# Row 24
$ws.Range("E24").Value = "  +4.42%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.90"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.72%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "93.78"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.07%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.23"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +4.39%  "

# Row 28
$ws.Range("E28").Value = "  +1.99%  "

# Row 29
$ws.Range("E29").Value = "  +1.09%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.47"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +2.95%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.21"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.33%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.32"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -1.68%  "

# Row 33
$ws.Range("E33").Value = "  +3.09%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "63.09"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -2.47%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.40"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +15.02%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.62"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +16.76%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "545.45"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -3.59%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.415"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +5.19%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "38.23"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.35%  "

# Row 40
$ws.Range("D40").Value = "0.0₃0809"
$ws.Range("E40").Value = "  +5.91%  "

# Row 41
$ws.Range("E41").Value = "  -0.11%  "

# Row 42
$ws.Range("D42").Value = "3.584.71"
$ws.Range("E42").Value = "  +9.83%  "

# Row 43
$ws.Range("E43").Value = "  +5.05%  "

# Row 44
$ws.Range("E44").Value = "  +3.69%  "

# Row 45
$ws.Range("E45").Value = "  +6.89%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.50"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.39%  "

# Row 47
$ws.Range("E47").Value = "  -0.85%  "

# Row 48
$ws.Range("E48").Value = "  +4.50%  "

# Row 49
$ws.Range("E49").Value = "  +2.79%  "

# Row 50
$ws.Range("E50").Value = "  +13.57%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.999"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.21%  "

